# Apply changes described by the diff:
# - D3, D5, D6, D7 phone numbers all updated to 3204886934
# - Active selection moved from C12 to F20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 3204886934
$ws.Range("D5").Value = 3204886934
$ws.Range("D6").Value = 3204886934
$ws.Range("D7").Value = 3204886934

# Move/set the active selection to F20 to match the sheetView selection in the diff
$ws.Range("F20").Select()
